$d = $word.ActiveDocument

# The document still stores "user content" markers as <w:fldSimple> fields
# (m:usercontent / m:endusercontent). When Word regenerates/touches such a
# field it normally "upgrades" it in place to the equivalent complex
# field-code construct: a run with fldChar begin, a run with the instrText,
# a run with fldChar separate and a run with fldChar end. Doing this once
# up-front (instead of leaving the <w:fldSimple> shortcut form) is what
# keeps the user content markers from getting lost across generations.
#
# Reproduce that upgrade for every simple field currently in the document,
# leaving everything else (surrounding paragraphs/runs/attributes) intact.

# Snapshot the field instructions first: InsertXML rebuilds the paragraph
# (and therefore the Fields/Paragraphs collections) as we go.
$instrList = New-Object System.Collections.ArrayList
foreach ($f in $d.Fields) {
    [void]$instrList.Add($f.Code.Text.Trim())
}

foreach ($instr in $instrList) {
    # Re-resolve the field (and its owning paragraph) on every iteration
    # since earlier replacements shift/rebuild ranges further in the body.
    $target = $null
    foreach ($f in $d.Fields) {
        if ($f.Code.Text.Trim() -eq $instr) {
            $target = $f
            break
        }
    }
    if ($target -eq $null) {
        continue
    }

    $paraIndex = 0
    $idx = 0
    foreach ($p in $d.Paragraphs) {
        $idx = $idx + 1
        if (($target.Code.Start -ge $p.Range.Start) -and ($target.Code.Start -lt $p.Range.End)) {
            $paraIndex = $idx
            break
        }
    }
    if ($paraIndex -eq 0) {
        continue
    }

    # Use the Paragraph object's own Range (covers the paragraph mark too)
    # so InsertXML fully replaces the <w:fldSimple> rather than inserting
    # the new runs alongside it.
    $rng = $d.Paragraphs.Item($paraIndex).Range

    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' +
           '<w:r><w:instrText>' + $instr + '</w:instrText></w:r>' +
           '<w:r><w:fldChar w:fldCharType="separate"/></w:r>' +
           '<w:r><w:fldChar w:fldCharType="end"/></w:r>' +
           '</w:p>'

    $rng.InsertXML($xml)
}
